$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.026.91'
$ws.Range("E2").Value = '  -5.29%  '
$ws.Range("D3").Value = '3.298.54'
$ws.Range("E3").Value = '  -5.99%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("E6").Value = '  -8.54%  '
$ws.Range("E7").Value = '  -5.21%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '3.289.26'
$ws.Range("E9").Value = '  -6.07%  '
$ws.Range("E10").Value = '  -4.93%  '
$ws.Range("E11").Value = '  -4.41%  '
$ws.Range("E12").Value = '  -5.46%  '
$ws.Range("E13").Value = '  -5.89%  '
$ws.Range("E14").Value = '  -6.38%  '
$ws.Range("D15").Value = '3.832.74'
$ws.Range("E15").Value = '  -5.83%  '
$ws.Range("E16").Value = '  -4.33%  '
$ws.Range("D17").Value = '3.314.53'
$ws.Range("E17").Value = '  -5.52%  '
$ws.Range("E18").Value = '  -5.70%  '
$ws.Range("E19").Value = '  -4.81%  '
$ws.Range("D20").Value = '63.072.48'
$ws.Range("E20").Value = '  -5.28%  '
$ws.Range("E21").Value = '  -4.44%  '
$ws.Range("E22").Value = '  -4.48%  '
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  -4.66%  '
$ws.Range("E26").Value = '  +4.94%  '
$ws.Range("E27").Value = '  -4.42%  '
$ws.Range("E28").Value = '  -7.82%  '
$ws.Range("E29").Value = '  -8.06%  '
$ws.Range("E30").Value = '  -5.38%  '
$ws.Range("E31").Value = '  -5.65%  '
$ws.Range("E32").Value = '  -6.05%  '
$ws.Range("E33").Value = '  -8.28%  '
$ws.Range("E34").Value = '  -6.29%  '
$ws.Range("E35").Value = '  -4.92%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -4.11%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E38").Value = '  -9.11%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").Value = '3.138.10'
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("E41").Value = '  -10.94%  '
$ws.Range("E42").Value = '  -6.51%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E44").Value = '  -3.15%  '
$ws.Range("E45").Value = '  -4.25%  '
$ws.Range("E46").Value = '  -4.90%  '
$ws.Range("E47").Value = '  -8.89%  '
$ws.Range("E49").Value = '  -5.62%  '
$ws.Range("E50").Value = '  -5.26%  '
$ws.Range("E51").Value = '  -7.06%  '

# Cells whose new value is a plain number-looking string; force Text
# so Excel stores them as text (matching the source data) instead of
# auto-converting to a numeric cell, then restore the default style.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '52.60'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.116'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '400.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '568.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.145'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '34.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.362'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.126'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.93'
$ws.Range("D51").Style = "Normal"
